$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "environment" row (row 6) NroSiniestro value
# (leading apostrophe keeps the existing text/quote-prefix cell style)
$ws.Range("E6").Value = "'0420172010449"

# Clear the extra claim-number rows (7 and 8), keeping formatting
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()

# Update the active selection to reflect the saved view state
$ws.Range("N6").Select()
